{"js": "const replacements = [\n  [\"4+43=\", \"58+0=\"],\n  [\"3+76=\", \"53+11=\"],\n  [\"87-55=\", \"59+25=\"],\n  [\"29+0=\", \"66-17=\"],\n  [\"87-37=\", \"93-16=\"],\n  [\"28+40=\", \"9+35=\"],\n  [\"11+78=\", \"88-13=\"],\n  [\"95-49=\", \"38+23=\"],\n  [\"43+9=\", \"20+22=\"],\n  [\"36-23=\", \"61-43=\"],\n  [\"43+50=\", \"95-47=\"],\n  [\"37+27=\", \"90-53=\"],\n  [\"71-42=\", \"28-11=\"],\n  [\"12+18=\", \"72-50=\"],\n  [\"90-11=\", \"48+19=\"],\n  [\"18+70=\", \"73-10=\"],\n  [\"40-16=\", \"92-51=\"],\n  [\"43+7=\", \"64+17=\"],\n  [\"1+62=\", \"49-7=\"],\n  [\"83-49=\", \"16+23=\"],\n  [\"58-54=\", \"87+10=\"],\n  [\"67-37=\", \"95-6=\"],\n  [\"82-74=\", \"39-23=\"],\n  [\"74-8=\", \"8+32=\"],\n  [\"26+2=\", \"62-4=\"],\n  [\"30+2=\", \"70-35=\"],\n  [\"68-13=\", \"31+28=\"],\n  [\"79-8=\", \"47+24=\"],\n  [\"29+69=\", \"48-29=\"],\n  [\"12+59=\", \"85-7=\"],\n  [\"0+51=\", \"83-71=\"],\n  [\"78-49=\", \"79-44=\"],\n  [\"37-6=\", \"86-35=\"],\n  [\"36+30=\", \"26+6=\"],\n  [\"78-14=\", \"10+28=\"],\n  [\"14+56=\", \"46+30=\"],\n  [\"66-51=\", \"34+17=\"],\n  [\"45+14=\", \"81-78=\"],\n  [\"66-45=\", \"81-24=\"],\n  [\"36+1=\", \"68-8=\"],\n  [\"27+44=\", \"24+69=\"],\n  [\"99-66=\", \"41+12=\"],\n  [\"10+89=\", \"88-6=\"],\n  [\"19+73=\", \"73-0=\"],\n  [\"23-22=\", \"99-41=\"],\n  [\"76-43=\", \"95-55=\"],\n  [\"48-25=\", \"5+86=\"],\n  [\"22+9=\", \"55-36=\"],\n  [\"89-79=\", \"29+50=\"],\n  [\"46-22=\", \"76-44=\"],\n  [\"88-75=\", \"67+1=\"],\n  [\"72+9=\", \"94-63=\"],\n  [\"81-15=\", \"55+32=\"],\n  [\"76-18=\", \"45-40=\"],\n  [\"29+20=\", \"87-60=\"],\n  [\"98-39=\", \"65-4=\"],\n  [\"88-64=\", \"66+30=\"],\n  [\"49+37=\", \"17+10=\"],\n  [\"23+59=\", \"81-74=\"],\n  [\"50+21=\", \"70-33=\"],\n  [\"36+40=\", \"22+42=\"],\n  [\"3+39=\", \"25+50=\"],\n  [\"8+59=\", \"17-10=\"],\n  [\"18+45=\", \"59-19=\"],\n  [\"33-5=\", \"87-1=\"],\n  [\"55-48=\", \"61-38=\"],\n  [\"97-22=\", \"62+35=\"],\n  [\"60+11=\", \"74-6=\"],\n  [\"99-1=\", \"47-14=\"],\n  [\"83-72=\", \"53-21=\"],\n  [\"33+60=\", \"51-41=\"],\n  [\"83-13=\", \"22+15=\"],\n  [\"31-11=\", \"99-53=\"],\n  [\"46+9=\", \"81+13=\"],\n  [\"42+46=\", \"94-2=\"],\n  [\"60+28=\", \"19-5=\"],\n  [\"54-47=\", \"42+57=\"],\n  [\"23+61=\", \"89+7=\"],\n  [\"64-63=\", \"51-30=\"],\n  [\"18-1=\", \"57-10=\"],\n  [\"17+28=\", \"20+76=\"],\n  [\"43-12=\", \"8-0=\"],\n  [\"8+58=\", \"25+0=\"],\n  [\"66+9=\", \"6+58=\"],\n  [\"81-1=\", \"3+25=\"],\n  [\"32-5=\", \"50-6=\"],\n  [\"1+63=\", \"40-23=\"],\n  [\"55-19=\", \"45+13=\"],\n  [\"18+36=\", \"3+31=\"],\n  [\"15+82=\", \"77+8=\"],\n  [\"46-34=\", \"72-42=\"],\n  [\"60-4=\", \"24+24=\"],\n  [\"62+21=\", \"14+85=\"],\n  [\"84-62=\", \"88-8=\"],\n  [\"39-20=\", \"63+36=\"],\n  [\"84-24=\", \"48+4=\"],\n  [\"26-9=\", \"72-41=\"],\n  [\"24+1=\", \"29+61=\"],\n  [\"83+1=\", \"63+18=\"],\n  [\"11+30=\", \"82+7=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"4+43=\", \"58+0=\"),\n    @(\"3+76=\", \"53+11=\"),\n    @(\"87-55=\", \"59+25=\"),\n    @(\"29+0=\", \"66-17=\"),\n    @(\"87-37=\", \"93-16=\"),\n    @(\"28+40=\", \"9+35=\"),\n    @(\"11+78=\", \"88-13=\"),\n    @(\"95-49=\", \"38+23=\"),\n    @(\"43+9=\", \"20+22=\"),\n    @(\"36-23=\", \"61-43=\"),\n    @(\"43+50=\", \"95-47=\"),\n    @(\"37+27=\", \"90-53=\"),\n    @(\"71-42=\", \"28-11=\"),\n    @(\"12+18=\", \"72-50=\"),\n    @(\"90-11=\", \"48+19=\"),\n    @(\"18+70=\", \"73-10=\"),\n    @(\"40-16=\", \"92-51=\"),\n    @(\"43+7=\", \"64+17=\"),\n    @(\"1+62=\", \"49-7=\"),\n    @(\"83-49=\", \"16+23=\"),\n    @(\"58-54=\", \"87+10=\"),\n    @(\"67-37=\", \"95-6=\"),\n    @(\"82-74=\", \"39-23=\"),\n    @(\"74-8=\", \"8+32=\"),\n    @(\"26+2=\", \"62-4=\"),\n    @(\"30+2=\", \"70-35=\"),\n    @(\"68-13=\", \"31+28=\"),\n    @(\"79-8=\", \"47+24=\"),\n    @(\"29+69=\", \"48-29=\"),\n    @(\"12+59=\", \"85-7=\"),\n    @(\"0+51=\", \"83-71=\"),\n    @(\"78-49=\", \"79-44=\"),\n    @(\"37-6=\", \"86-35=\"),\n    @(\"36+30=\", \"26+6=\"),\n    @(\"78-14=\", \"10+28=\"),\n    @(\"14+56=\", \"46+30=\"),\n    @(\"66-51=\", \"34+17=\"),\n    @(\"45+14=\", \"81-78=\"),\n    @(\"66-45=\", \"81-24=\"),\n    @(\"36+1=\", \"68-8=\"),\n    @(\"27+44=\", \"24+69=\"),\n    @(\"99-66=\", \"41+12=\"),\n    @(\"10+89=\", \"88-6=\"),\n    @(\"19+73=\", \"73-0=\"),\n    @(\"23-22=\", \"99-41=\"),\n    @(\"76-43=\", \"95-55=\"),\n    @(\"48-25=\", \"5+86=\"),\n    @(\"22+9=\", \"55-36=\"),\n    @(\"89-79=\", \"29+50=\"),\n    @(\"46-22=\", \"76-44=\"),\n    @(\"88-75=\", \"67+1=\"),\n    @(\"72+9=\", \"94-63=\"),\n    @(\"81-15=\", \"55+32=\"),\n    @(\"76-18=\", \"45-40=\"),\n    @(\"29+20=\", \"87-60=\"),\n    @(\"98-39=\", \"65-4=\"),\n    @(\"88-64=\", \"66+30=\"),\n    @(\"49+37=\", \"17+10=\"),\n    @(\"23+59=\", \"81-74=\"),\n    @(\"50+21=\", \"70-33=\"),\n    @(\"36+40=\", \"22+42=\"),\n    @(\"3+39=\", \"25+50=\"),\n    @(\"8+59=\", \"17-10=\"),\n    @(\"18+45=\", \"59-19=\"),\n    @(\"33-5=\", \"87-1=\"),\n    @(\"55-48=\", \"61-38=\"),\n    @(\"97-22=\", \"62+35=\"),\n    @(\"60+11=\", \"74-6=\"),\n    @(\"99-1=\", \"47-14=\"),\n    @(\"83-72=\", \"53-21=\"),\n    @(\"33+60=\", \"51-41=\"),\n    @(\"83-13=\", \"22+15=\"),\n    @(\"31-11=\", \"99-53=\"),\n    @(\"46+9=\", \"81+13=\"),\n    @(\"42+46=\", \"94-2=\"),\n    @(\"60+28=\", \"19-5=\"),\n    @(\"54-47=\", \"42+57=\"),\n    @(\"23+61=\", \"89+7=\"),\n    @(\"64-63=\", \"51-30=\"),\n    @(\"18-1=\", \"57-10=\"),\n    @(\"17+28=\", \"20+76=\"),\n    @(\"43-12=\", \"8-0=\"),\n    @(\"8+58=\", \"25+0=\"),\n    @(\"66+9=\", \"6+58=\"),\n    @(\"81-1=\", \"3+25=\"),\n    @(\"32-5=\", \"50-6=\"),\n    @(\"1+63=\", \"40-23=\"),\n    @(\"55-19=\", \"45+13=\"),\n    @(\"18+36=\", \"3+31=\"),\n    @(\"15+82=\", \"77+8=\"),\n    @(\"46-34=\", \"72-42=\"),\n    @(\"60-4=\", \"24+24=\"),\n    @(\"62+21=\", \"14+85=\"),\n    @(\"84-62=\", \"88-8=\"),\n    @(\"39-20=\", \"63+36=\"),\n    @(\"84-24=\", \"48+4=\"),\n    @(\"26-9=\", \"72-41=\"),\n    @(\"24+1=\", \"29+61=\"),\n    @(\"83+1=\", \"63+18=\"),\n    @(\"11+30=\", \"82+7=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
